# Gnai2-Adora1.xlsx - "update scripts wuth new tpm"
#
# The NATMI TPM recompute shifted every LR-pair row: the "Target cluster"
# (col D) rotates ECs->MuSCs / MuSCs->Resolving-Mac, and every downstream
# expression / specificity statistic (cols G:T) is refreshed with the new
# TPM-derived numbers. Columns A, B, C, E, F are untouched by the update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> ("Target cluster" string, [G,H,I,J,K,L,M,N,O,P,Q,R,S,T] or $null to leave as-is)
$rows = @(
    @{ Row = 2;  Target = "MuSCs";         Values = @(201.4397426666667, 604.3192280000001, 0.4833500233086392, 0.4833500233086393, 3, 1,                    0.782436,            2.347308,  0.9918763017460563, 0.9918763017460563, 157.613706493136,   1418.523358438224, 0.4794234335682432,  0.4794234335682433) },
    @{ Row = 3;  Target = "Resolving-Mac"; Values = @(201.4397426666667, 604.3192280000001, 0.4833500233086392, 0.4833500233086393, 2, 0.6666666666666666,   0.006408333333333333, 0.019225, 0.008123698253943637, 0.008123698253943637, 1.290893017588889, 11.6180371583,      0.003926589740396009, 0.00392658974039601) },
    @{ Row = 4;  Target = "MuSCs";         Values = @($null, $null,      0.1569674599353791, 0.1569674599353792, 3, 1,                    0.782436,            2.347308,  0.9918763017460563, 0.9918763017460563, 51.184901140328,    460.664110262952,  0.1556923036551761,  0.1556923036551762) },
    @{ Row = 5;  Target = "Resolving-Mac"; Values = @($null, $null,      0.1569674599353791, 0.1569674599353792, 2, 0.6666666666666666,   0.006408333333333333, 0.019225, 0.008123698253943637, 0.008123698253943637, 0.4192162785722222, 3.77294650715,     0.001275156280203007, 0.001275156280203008) },
    @{ Row = 6;  Target = "MuSCs";         Values = @(60.43484133333334, 181.304524,         0.1450120099461104, 0.1450120099461104, 3, 1,                    0.782436,            2.347308,  0.9918763017460563, 0.9918763017460563, 47.286395513488,    425.577559621392,  0.1438339761341103,  0.1438339761341103) },
    @{ Row = 7;  Target = "Resolving-Mac"; Values = @(60.43484133333334, 181.304524,         0.1450120099461104, 0.1450120099461104, 2, 0.6666666666666666,   0.006408333333333333, 0.019225, 0.008123698253943637, 0.008123698253943637, 0.3872866082111111, 3.4855794739,      0.001178033812000074, 0.001178033812000074) },
    @{ Row = 8;  Target = "MuSCs";         Values = @(89.46554166666668, 268.396625,         0.2146705068098712, 0.2146705068098712, 3, 1,                    0.782436,            2.347308,  0.9918763017460563, 0.9918763017460563, 70.00106055950002,  630.0095450355001, 0.2129265883885267,  0.2129265883885267) },
    @{ Row = 9;  Target = "Resolving-Mac"; Values = @(89.46554166666668, 268.396625,         0.2146705068098712, 0.2146705068098712, 2, 0.6666666666666666,   0.006408333333333333, 0.019225, 0.008123698253943637, 0.008123698253943637, 0.5733250128472223, 5.159925115625001, 0.001743918421344546, 0.001743918421344547) }
)

# Columns G..T, in order, matching the Values arrays above.
$cols = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $rows) {
    $r = $entry.Row

    # Col D: "Target cluster"
    $ws.Range("D$r").Value = $entry.Target

    # Cols G:T
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $v = $entry.Values[$i]
        if ($null -ne $v) {
            $ws.Range("$($cols[$i])$r").Value = $v
        }
    }
}
